$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

# Row 96
$ws.Range("A96").Value = 130964396
$ws.Range("B96").Value = 98931
$ws.Range("E96").Value = 219790
$ws.Range("Q96").Value = 509365
$ws.Range("R96").Value = 6718883
$ws.Range("D96").Value = 'LC'
$ws.Range("F96").Value = 'Fläcknycklar'
$ws.Range("G96").Value = 'Dactylorhiza maculata'
$ws.Range("H96").Value = '(L.) Soó'
$ws.Range("AC96").Value = 'Betydande förekomster . inventering åt vasa vind'
$ws.Range("Y96").Formula = '="2025-07-03"'
$ws.Range("Y96").Copy()
$ws.Range("Y96").PasteSpecial(-4163)
$ws.Range("AA96").Formula = '="2025-07-03"'
$ws.Range("AA96").Copy()
$ws.Range("AA96").PasteSpecial(-4163)

# Row 97
$ws.Range("A97").Value = 130964529
$ws.Range("B97").Value = 79244
$ws.Range("E97").Value = 6425
$ws.Range("Q97").Value = 509588
$ws.Range("R97").Value = 6719099
$ws.Range("D97").Value = 'NT'
$ws.Range("F97").Value = 'Garnlav'
$ws.Range("G97").Value = 'Alectoria sarmentosa'
$ws.Range("H97").Value = '(Ach.) Ach.'
$ws.Range("AC97").Value = 'Enstaka . inventering åt vasa vind'
$ws.Range("Y97").Formula = '="2025-07-02"'
$ws.Range("Y97").Copy()
$ws.Range("Y97").PasteSpecial(-4163)
$ws.Range("AA97").Formula = '="2025-07-02"'
$ws.Range("AA97").Copy()
$ws.Range("AA97").PasteSpecial(-4163)

# Row 98
$ws.Range("A98").Value = 130964573
$ws.Range("B98").Value = 79244
$ws.Range("E98").Value = 6425
$ws.Range("Q98").Value = 509515
$ws.Range("R98").Value = 6719063
$ws.Range("D98").Value = 'NT'
$ws.Range("F98").Value = 'Garnlav'
$ws.Range("G98").Value = 'Alectoria sarmentosa'
$ws.Range("H98").Value = '(Ach.) Ach.'
$ws.Range("AC98").Value = 'Måttliga förekomster . inventering åt vasa vind'
$ws.Range("Y98").Formula = '="2025-07-02"'
$ws.Range("Y98").Copy()
$ws.Range("Y98").PasteSpecial(-4163)
$ws.Range("AA98").Formula = '="2025-07-02"'
$ws.Range("AA98").Copy()
$ws.Range("AA98").PasteSpecial(-4163)

# Row 102
$ws.Range("A102").Value = 130964526
$ws.Range("B102").Value = 79244
$ws.Range("E102").Value = 6425
$ws.Range("Q102").Value = 509610
$ws.Range("R102").Value = 6719050
$ws.Range("D102").Value = 'NT'
$ws.Range("F102").Value = 'Garnlav'
$ws.Range("G102").Value = 'Alectoria sarmentosa'
$ws.Range("H102").Value = '(Ach.) Ach.'
$ws.Range("AC102").Value = 'Enstaka . inventering åt vasa vind'
$ws.Range("Y102").Formula = '="2025-07-02"'
$ws.Range("Y102").Copy()
$ws.Range("Y102").PasteSpecial(-4163)
$ws.Range("AA102").Formula = '="2025-07-02"'
$ws.Range("AA102").Copy()
$ws.Range("AA102").PasteSpecial(-4163)

# Row 103
$ws.Range("A103").Value = 130964547
$ws.Range("B103").Value = 57881
$ws.Range("E103").Value = 100049
$ws.Range("Q103").Value = 509495
$ws.Range("R103").Value = 6718877
$ws.Range("D103").Value = 'NT'
$ws.Range("F103").Value = 'Spillkråka'
$ws.Range("G103").Value = 'Dryocopus martius'
$ws.Range("H103").Value = '(Linnaeus, 1758)'
$ws.Range("AC103").Value = 'Födosökspår . inventering åt vasa vind'
$ws.Range("Y103").Formula = '="2025-07-02"'
$ws.Range("Y103").Copy()
$ws.Range("Y103").PasteSpecial(-4163)
$ws.Range("AA103").Formula = '="2025-07-02"'
$ws.Range("AA103").Copy()
$ws.Range("AA103").PasteSpecial(-4163)

# Row 105
$ws.Range("A105").Value = 130964541
$ws.Range("B105").Value = 91809
$ws.Range("E105").Value = 1202
$ws.Range("Q105").Value = 509703
$ws.Range("R105").Value = 6719018
$ws.Range("D105").Value = 'NT'
$ws.Range("F105").Value = 'Ullticka'
$ws.Range("G105").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H105").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("AC105").Value = 'Enstaka . inventering åt vasa vind'
$ws.Range("Y105").Formula = '="2025-07-02"'
$ws.Range("Y105").Copy()
$ws.Range("Y105").PasteSpecial(-4163)
$ws.Range("AA105").Formula = '="2025-07-02"'
$ws.Range("AA105").Copy()
$ws.Range("AA105").PasteSpecial(-4163)

# Row 106
$ws.Range("A106").Value = 130964537
$ws.Range("B106").Value = 79244
$ws.Range("E106").Value = 6425
$ws.Range("Q106").Value = 509822
$ws.Range("R106").Value = 6718960
$ws.Range("D106").Value = 'NT'
$ws.Range("F106").Value = 'Garnlav'
$ws.Range("G106").Value = 'Alectoria sarmentosa'
$ws.Range("H106").Value = '(Ach.) Ach.'
$ws.Range("AC106").Value = 'Rikligt . inventering åt vasa vind'
$ws.Range("Y106").Formula = '="2025-07-02"'
$ws.Range("Y106").Copy()
$ws.Range("Y106").PasteSpecial(-4163)
$ws.Range("AA106").Formula = '="2025-07-02"'
$ws.Range("AA106").Copy()
$ws.Range("AA106").PasteSpecial(-4163)

# Row 108
$ws.Range("A108").Value = 130964574
$ws.Range("B108").Value = 79244
$ws.Range("E108").Value = 6425
$ws.Range("Q108").Value = 509667
$ws.Range("R108").Value = 6719184
$ws.Range("D108").Value = 'NT'
$ws.Range("F108").Value = 'Garnlav'
$ws.Range("G108").Value = 'Alectoria sarmentosa'
$ws.Range("H108").Value = '(Ach.) Ach.'
$ws.Range("AC108").Value = 'Måttliga förekomster . inventering åt vasa vind'
$ws.Range("Y108").Formula = '="2025-07-02"'
$ws.Range("Y108").Copy()
$ws.Range("Y108").PasteSpecial(-4163)
$ws.Range("AA108").Formula = '="2025-07-02"'
$ws.Range("AA108").Copy()
$ws.Range("AA108").PasteSpecial(-4163)

# Row 109
$ws.Range("A109").Value = 130964544
$ws.Range("B109").Value = 57073
$ws.Range("E109").Value = 100138
$ws.Range("Q109").Value = 509543
$ws.Range("R109").Value = 6718926
$ws.Range("D109").Value = 'LC'
$ws.Range("F109").Value = 'Tjäder'
$ws.Range("G109").Value = 'Tetrao urogallus'
$ws.Range("H109").Value = 'Linnaeus, 1758'
$ws.Range("AC109").Value = 'Spillning . inventering åt vasa vind'
$ws.Range("Y109").Formula = '="2025-07-02"'
$ws.Range("Y109").Copy()
$ws.Range("Y109").PasteSpecial(-4163)
$ws.Range("AA109").Formula = '="2025-07-02"'
$ws.Range("AA109").Copy()
$ws.Range("AA109").PasteSpecial(-4163)

# Row 111
$ws.Range("A111").Value = 130964641
$ws.Range("B111").Value = 98931
$ws.Range("E111").Value = 219790
$ws.Range("Q111").Value = 509932
$ws.Range("R111").Value = 6719045
$ws.Range("D111").Value = 'LC'
$ws.Range("F111").Value = 'Fläcknycklar'
$ws.Range("G111").Value = 'Dactylorhiza maculata'
$ws.Range("H111").Value = '(L.) Soó'
$ws.Range("AC111").Value = 'Måttlig förekomst . inventering åt vasa vind'
$ws.Range("Y111").Formula = '="2025-07-02"'
$ws.Range("Y111").Copy()
$ws.Range("Y111").PasteSpecial(-4163)
$ws.Range("AA111").Formula = '="2025-07-02"'
$ws.Range("AA111").Copy()
$ws.Range("AA111").PasteSpecial(-4163)

# Row 112
$ws.Range("A112").Value = 130964650
$ws.Range("B112").Value = 92268
$ws.Range("E112").Value = 1209
$ws.Range("Q112").Value = 509694
$ws.Range("R112").Value = 6718936
$ws.Range("D112").Value = 'VU'
$ws.Range("F112").Value = 'Rynkskinn'
$ws.Range("G112").Value = 'Hermanssonia centrifuga'
$ws.Range("H112").Value = '(P. Karst.) Zmitr.'
$ws.Range("AC112").Value = 'Måttliga förekomster . inventering åt vasa vind'
$ws.Range("Y112").Formula = '="2025-07-02"'
$ws.Range("Y112").Copy()
$ws.Range("Y112").PasteSpecial(-4163)
$ws.Range("AA112").Formula = '="2025-07-02"'
$ws.Range("AA112").Copy()
$ws.Range("AA112").PasteSpecial(-4163)

# Row 113
$ws.Range("A113").Value = 130964533
$ws.Range("B113").Value = 79244
$ws.Range("E113").Value = 6425
$ws.Range("Q113").Value = 509984
$ws.Range("R113").Value = 6719028
$ws.Range("D113").Value = 'NT'
$ws.Range("F113").Value = 'Garnlav'
$ws.Range("G113").Value = 'Alectoria sarmentosa'
$ws.Range("H113").Value = '(Ach.) Ach.'
$ws.Range("AC113").Value = 'Rikligt . inventering åt vasa vind'
$ws.Range("Y113").Formula = '="2025-07-02"'
$ws.Range("Y113").Copy()
$ws.Range("Y113").PasteSpecial(-4163)
$ws.Range("AA113").Formula = '="2025-07-02"'
$ws.Range("AA113").Copy()
$ws.Range("AA113").PasteSpecial(-4163)

# Row 114
$ws.Range("A114").Value = 130964645
$ws.Range("B114").Value = 99037
$ws.Range("E114").Value = 221952
$ws.Range("Q114").Value = 509804
$ws.Range("R114").Value = 6719024
$ws.Range("D114").Value = 'LC'
$ws.Range("F114").Value = 'Spindelblomster'
$ws.Range("G114").Value = 'Neottia cordata'
$ws.Range("H114").Value = '(L.) Rich.'
$ws.Range("AC114").Value = 'Måttliga förekomster . inventering åt vasa vind'
$ws.Range("Y114").Formula = '="2025-07-02"'
$ws.Range("Y114").Copy()
$ws.Range("Y114").PasteSpecial(-4163)
$ws.Range("AA114").Formula = '="2025-07-02"'
$ws.Range("AA114").Copy()
$ws.Range("AA114").PasteSpecial(-4163)

# Row 119
$ws.Range("A119").Value = 130964649
$ws.Range("B119").Value = 98931
$ws.Range("E119").Value = 219790
$ws.Range("Q119").Value = 509705
$ws.Range("R119").Value = 6718923
$ws.Range("D119").Value = 'LC'
$ws.Range("F119").Value = 'Fläcknycklar'
$ws.Range("G119").Value = 'Dactylorhiza maculata'
$ws.Range("H119").Value = '(L.) Soó'
$ws.Range("AC119").Value = 'Måttlig förekomst . inventering åt vasa vind'
$ws.Range("Y119").Formula = '="2025-07-02"'
$ws.Range("Y119").Copy()
$ws.Range("Y119").PasteSpecial(-4163)
$ws.Range("AA119").Formula = '="2025-07-02"'
$ws.Range("AA119").Copy()
$ws.Range("AA119").PasteSpecial(-4163)

# Row 120
$ws.Range("A120").Value = 130964648
$ws.Range("B120").Value = 92268
$ws.Range("E120").Value = 1209
$ws.Range("Q120").Value = 509744
$ws.Range("R120").Value = 6718982
$ws.Range("D120").Value = 'VU'
$ws.Range("F120").Value = 'Rynkskinn'
$ws.Range("G120").Value = 'Hermanssonia centrifuga'
$ws.Range("H120").Value = '(P. Karst.) Zmitr.'
$ws.Range("AC120").Value = 'Måttliga förekomster . inventering åt vasa vind'
$ws.Range("Y120").Formula = '="2025-07-02"'
$ws.Range("Y120").Copy()
$ws.Range("Y120").PasteSpecial(-4163)
$ws.Range("AA120").Formula = '="2025-07-02"'
$ws.Range("AA120").Copy()
$ws.Range("AA120").PasteSpecial(-4163)

# Row 122
$ws.Range("A122").Value = 130964644
$ws.Range("B122").Value = 98918
$ws.Range("E122").Value = 220093
$ws.Range("Q122").Value = 509801
$ws.Range("R122").Value = 6719017
$ws.Range("D122").Value = 'LC'
$ws.Range("F122").Value = 'Korallrot'
$ws.Range("G122").Value = 'Corallorhiza trifida'
$ws.Range("H122").Value = 'Châtel.'
$ws.Range("AC122").Value = 'Sparsamma förekomster . inventering åt vasa vind'
$ws.Range("Y122").Formula = '="2025-07-02"'
$ws.Range("Y122").Copy()
$ws.Range("Y122").PasteSpecial(-4163)
$ws.Range("AA122").Formula = '="2025-07-02"'
$ws.Range("AA122").Copy()
$ws.Range("AA122").PasteSpecial(-4163)

# Row 124
$ws.Range("A124").Value = 130964542
$ws.Range("B124").Value = 57073
$ws.Range("E124").Value = 100138
$ws.Range("Q124").Value = 509635
$ws.Range("R124").Value = 6718941
$ws.Range("D124").Value = 'LC'
$ws.Range("F124").Value = 'Tjäder'
$ws.Range("G124").Value = 'Tetrao urogallus'
$ws.Range("H124").Value = 'Linnaeus, 1758'
$ws.Range("AC124").Value = 'Spillning . inventering åt vasa vind'
$ws.Range("Y124").Formula = '="2025-07-02"'
$ws.Range("Y124").Copy()
$ws.Range("Y124").PasteSpecial(-4163)
$ws.Range("AA124").Formula = '="2025-07-02"'
$ws.Range("AA124").Copy()
$ws.Range("AA124").PasteSpecial(-4163)

# Row 125
$ws.Range("A125").Value = 130964390
$ws.Range("B125").Value = 99014
$ws.Range("E125").Value = 220787
$ws.Range("Q125").Value = 509475
$ws.Range("R125").Value = 6718881
$ws.Range("D125").Value = 'VU'
$ws.Range("F125").Value = 'Knärot'
$ws.Range("G125").Value = 'Goodyera repens'
$ws.Range("H125").Value = '(L.) R. Br.'
$ws.Range("AC125").Value = 'Måttliga förekomster, Ca 10-15 plantor . inventering åt vasa vind'
$ws.Range("Y125").Formula = '="2025-07-03"'
$ws.Range("Y125").Copy()
$ws.Range("Y125").PasteSpecial(-4163)
$ws.Range("AA125").Formula = '="2025-07-03"'
$ws.Range("AA125").Copy()
$ws.Range("AA125").PasteSpecial(-4163)

# Row 126
$ws.Range("A126").Value = 130964643
$ws.Range("B126").Value = 98931
$ws.Range("E126").Value = 219790
$ws.Range("Q126").Value = 509829
$ws.Range("R126").Value = 6719000
$ws.Range("D126").Value = 'LC'
$ws.Range("F126").Value = 'Fläcknycklar'
$ws.Range("G126").Value = 'Dactylorhiza maculata'
$ws.Range("H126").Value = '(L.) Soó'
$ws.Range("AC126").Value = 'Måttlig förekomst . inventering åt vasa vind'
$ws.Range("Y126").Formula = '="2025-07-02"'
$ws.Range("Y126").Copy()
$ws.Range("Y126").PasteSpecial(-4163)
$ws.Range("AA126").Formula = '="2025-07-02"'
$ws.Range("AA126").Copy()
$ws.Range("AA126").PasteSpecial(-4163)

# Row 127
$ws.Range("A127").Value = 130964546
$ws.Range("B127").Value = 92504
$ws.Range("E127").Value = 898
$ws.Range("Q127").Value = 509515
$ws.Range("R127").Value = 6718886
$ws.Range("D127").Value = 'VU'
$ws.Range("F127").Value = 'Blackticka'
$ws.Range("G127").Value = 'Steccherinum collabens'
$ws.Range("H127").Value = '(Fr.) Vesterholt'
$ws.Range("AC127").Value = 'Betydande förekomst . inventering åt vasa vind'
$ws.Range("Y127").Formula = '="2025-07-02"'
$ws.Range("Y127").Copy()
$ws.Range("Y127").PasteSpecial(-4163)
$ws.Range("AA127").Formula = '="2025-07-02"'
$ws.Range("AA127").Copy()
$ws.Range("AA127").PasteSpecial(-4163)

# Row 128
$ws.Range("A128").Value = 130964538
$ws.Range("B128").Value = 79244
$ws.Range("E128").Value = 6425
$ws.Range("Q128").Value = 509875
$ws.Range("R128").Value = 6719025
$ws.Range("D128").Value = 'NT'
$ws.Range("F128").Value = 'Garnlav'
$ws.Range("G128").Value = 'Alectoria sarmentosa'
$ws.Range("H128").Value = '(Ach.) Ach.'
$ws.Range("AC128").Value = 'Enstaka . inventering åt vasa vind'
$ws.Range("Y128").Formula = '="2025-07-02"'
$ws.Range("Y128").Copy()
$ws.Range("Y128").PasteSpecial(-4163)
$ws.Range("AA128").Formula = '="2025-07-02"'
$ws.Range("AA128").Copy()
$ws.Range("AA128").PasteSpecial(-4163)
